$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the rule table: header cell C1 "Investor Approval" -> "Auto Doc Gen"
$ws.Range("C1").Value = "Auto Doc Gen"

# 2. Insert a new blank row above the old row 3 (the "RuleTable InvestorApproval"
#    banner row), pushing everything below down by one row.
$ws.Rows("2:2").Insert()

# 3. Format the newly inserted row 2 (A2:D2) to match the orange accent-bar
#    look used elsewhere in the sheet (font colour + fill copied from the
#    existing banner style, fill recoloured to the orange accent).
$ws.Range("B1:D1").Copy()
$ws.Range("B2:D2").PasteSpecial(-4122)
$ws.Range("B2:D2").Interior.ThemeColor = 6
$ws.Range("B2:D2").Interior.TintAndShade = 0

# 4. Rename the rule table banner (now on row 4): "RuleTable InvestorApproval" -> "RuleTable AutoDocGen"
$ws.Range("B4").Value = "RuleTable AutoDocGen"

# 5. Restore the selection / view state to match the edited workbook.
$ws.Range("B5").Select()
